# Refresh cached market-board averages / leve profit figures across several
# class sheets (scheduled market-data sync). Each block below updates the
# currentAveragePrice* (H/I/J), LevePrice* (K/L) and LeveProfit* (M/N)
# columns for specific leve rows; some rows gain/lose the HQ profit (N) cell
# when HQ pricing becomes (un)available.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1709.6
$ws.Range("I19").Value = 1674.6666
$ws.Range("K19").Value = 1674.6666
$ws.Range("M19").Value = -1499.6666

$ws.Range("H43").Value = 798.3333
$ws.Range("I43").Value = 798
$ws.Range("K43").Value = 798
$ws.Range("M43").Value = -729

$ws.Range("H51").Value = 5499.5
$ws.Range("I51").Value = 3999
$ws.Range("J51").Value = 7000
$ws.Range("K51").Value = 3999
$ws.Range("L51").Value = 7000
$ws.Range("M51").Value = -3515
$ws.Range("N51").Value = -7968

$ws.Range("H70").Value = 3764.9333
$ws.Range("I70").Value = 2799
$ws.Range("J70").Value = 4006.4167
$ws.Range("K70").Value = 8397
$ws.Range("L70").Value = 12019.2501
$ws.Range("M70").Value = -8127
$ws.Range("N70").Value = -12559.2501

$ws.Range("H73").Value = 3764.9333
$ws.Range("I73").Value = 2799
$ws.Range("J73").Value = 4006.4167
$ws.Range("K73").Value = 8397
$ws.Range("L73").Value = 12019.2501
$ws.Range("M73").Value = -7461
$ws.Range("N73").Value = -13891.2501

$ws.Range("H111").Value = 409.33334
$ws.Range("I111").Value = 409.33334
$ws.Range("K111").Value = 1228.00002
$ws.Range("M111").Value = 1838.99998

$ws.Range("H137").Value = 1780.0385
$ws.Range("I137").Value = 1568.25
$ws.Range("J137").Value = 2486
$ws.Range("K137").Value = 4704.75
$ws.Range("L137").Value = 7458
$ws.Range("M137").Value = -2154.75
$ws.Range("N137").Value = -12558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3121.8572
$ws.Range("I61").Value = 3121.8572
$ws.Range("K61").Value = 3121.8572
$ws.Range("M61").Value = -2909.8572

$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

$ws.Range("H102").Value = 1209.6
$ws.Range("I102").Value = 917.6667
$ws.Range("J102").Value = 2377.3333
$ws.Range("K102").Value = 917.6667
$ws.Range("L102").Value = 2377.3333
$ws.Range("M102").Value = 704.3333
$ws.Range("N102").Value = -5621.3333

$ws.Range("H110").Value = 587.875
$ws.Range("I110").Value = 451.33334
$ws.Range("J110").Value = 997.5
$ws.Range("K110").Value = 451.33334
$ws.Range("L110").Value = 997.5
$ws.Range("M110").Value = 1593.66666
$ws.Range("N110").Value = -5087.5

$ws.Range("H132").Value = 1878.3636
$ws.Range("I132").Value = 1878.3636
$ws.Range("K132").Value = 5635.0908
$ws.Range("M132").Value = -3105.0908

$ws.Range("H136").Value = 3121.8572
$ws.Range("I136").Value = 3121.8572
$ws.Range("K136").Value = 9365.571599999999
$ws.Range("M136").Value = -6815.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3786.625
$ws.Range("J105").Value = 3199.5
$ws.Range("L105").Value = 3199.5
$ws.Range("N105").Value = -6693.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 925.1429000000001
$ws.Range("J16").Value = 950
$ws.Range("L16").Value = 950
$ws.Range("N16").Value = -1524

$ws.Range("H86").Value = 13047.286
$ws.Range("I86").Value = 13110.333
$ws.Range("J86").Value = 13000
$ws.Range("K86").Value = 13110.333
$ws.Range("L86").Value = 13000
$ws.Range("M86").Value = -11987.333
$ws.Range("N86").Value = -15246

$ws.Range("H89").Value = 13047.286
$ws.Range("I89").Value = 13110.333
$ws.Range("J89").Value = 13000
$ws.Range("K89").Value = 65551.66500000001
$ws.Range("L89").Value = 65000
$ws.Range("M89").Value = -59935.66500000001
$ws.Range("N89").Value = -76232

$ws.Range("H105").Value = 1167
$ws.Range("I105").Value = 822.6667
$ws.Range("J105").Value = 2200
$ws.Range("K105").Value = 822.6667
$ws.Range("L105").Value = 2200
$ws.Range("M105").Value = 924.3333
$ws.Range("N105").Value = -5694

$ws.Range("H107").Value = 682.8
$ws.Range("I107").Value = 249.5
$ws.Range("J107").Value = 971.6667
$ws.Range("K107").Value = 249.5
$ws.Range("L107").Value = 971.6667
$ws.Range("M107").Value = 1670.5
$ws.Range("N107").Value = -4811.6667

$ws.Range("H113").Value = 925.1429000000001
$ws.Range("J113").Value = 950
$ws.Range("L113").Value = 950
$ws.Range("N113").Value = -5290

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 99751.5
$ws.Range("J37").Value = 99751.5
$ws.Range("L37").Value = 299254.5
$ws.Range("N37").Value = -299478.5

$ws.Range("H68").Value = 1849.75
$ws.Range("I68").Value = 1799.6666
$ws.Range("K68").Value = 5398.9998
$ws.Range("M68").Value = -4587.9998

$ws.Range("H71").Value = 1849.75
$ws.Range("I71").Value = 1799.6666
$ws.Range("K71").Value = 16196.9994
$ws.Range("M71").Value = -12140.9994

$ws.Range("H86").Value = 986
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 986
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2485.1875
$ws.Range("I132").Value = 2533.3572
$ws.Range("K132").Value = 7600.071599999999
$ws.Range("M132").Value = -5070.071599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2933
$ws.Range("I40").Value = 2933
$ws.Range("K40").Value = 2933
$ws.Range("M40").Value = -2797

$ws.Range("H61").Value = 5747.3335
$ws.Range("I61").Value = 5500
$ws.Range("K61").Value = 5500
$ws.Range("M61").Value = -5298

$ws.Range("H68").Value = 27538.75
$ws.Range("I68").Value = 1703
$ws.Range("J68").Value = 53374.5
$ws.Range("K68").Value = 1703
$ws.Range("L68").Value = 53374.5
$ws.Range("M68").Value = -954
$ws.Range("N68").Value = -54872.5

$ws.Range("H71").Value = 27538.75
$ws.Range("I71").Value = 1703
$ws.Range("J71").Value = 53374.5
$ws.Range("K71").Value = 8515
$ws.Range("L71").Value = 266872.5
$ws.Range("M71").Value = -4771
$ws.Range("N71").Value = -274360.5

$ws.Range("H113").Value = 5747.3335
$ws.Range("I113").Value = 5500
$ws.Range("K113").Value = 5500
$ws.Range("M113").Value = -3330

$ws.Range("H122").Value = 3502.7144
$ws.Range("I122").Value = 3503.1667
$ws.Range("K122").Value = 10509.5001
$ws.Range("M122").Value = -8059.500100000001

$ws.Range("H141").Value = 59500
$ws.Range("J141").Value = 59500
$ws.Range("L141").Value = 59500
$ws.Range("N141").Value = -69860
